$wb = $excel.ActiveWorkbook

# ---- Sheet0 (sheet1.xml): add rows 78-81 ----
$ws0 = $wb.Worksheets.Item("Sheet0")
$ws0.Range("A77").Copy()
$ws0.Range("A78:A81").PasteSpecial(-4122)

$ws0.Range("A78").Value = 76
$ws0.Range("B78").Value = 85
$ws0.Range("C78").Value = 'C0045392F'
$ws0.Range("D78").Value = '2020-11-10'
$ws0.Range("E78").Value = '/Users/Keshab/Desktop/fewImages/1.jpg'

$ws0.Range("A79").Value = 77
$ws0.Range("B79").Value = 86
$ws0.Range("C79").Value = 'C0045392F'
$ws0.Range("D79").Value = '2020-11-10'
$ws0.Range("E79").Value = '/Users/Keshab/Desktop/fewImages/2.jpg'

$ws0.Range("A80").Value = 78
$ws0.Range("B80").Value = 87
$ws0.Range("C80").Value = 'C0601382F'
$ws0.Range("D80").Value = '2020-11-10'
$ws0.Range("E80").Value = '/Users/Keshab/Desktop/fewImages/3.jpg'

$ws0.Range("A81").Value = 79
$ws0.Range("B81").Value = 88
$ws0.Range("C81").Value = 'C0601382F'
$ws0.Range("D81").Value = '2020-11-10'
$ws0.Range("E81").Value = '/Users/Keshab/Desktop/fewImages/4.jpg'

# ---- Sheet1 (sheet2.xml): add rows 78-81 ----
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A77").Copy()
$ws1.Range("A78:A81").PasteSpecial(-4122)

$ws1.Range("A78").Value = 76
$ws1.Range("B78").Value = 85
$ws1.Range("C78").Value = ' JULIAN A. STEYERMARK , COLLECTOR'
$ws1.Range("D78").Value = ' May 9 , 198 %'
$ws1.Range("E78").Value = ' Falling Springs Duo . near'
$ws1.Range("F78").Value = ' No. 1671'
$ws1.Range("G78").Value = ' cystopteris alata'
$ws1.Range("H78").Value = ' FLORA OF MISSOURI Jllinois'

$ws1.Range("A79").Value = 77
$ws1.Range("B79").Value = 86
$ws1.Range("C79").Value = ' JULIAN A. STEYERMARK , COLLECTOR'
$ws1.Range("D79").Value = ' May 9 , 1984'
$ws1.Range("E79").Value = ' ( L. ) Bernh . Falling Springs near .'
$ws1.Range("F79").Value = ' No.1677'
$ws1.Range("G79").Value = ' Cystopteris pellucida'
$ws1.Range("H79").Value = ' FLORA OF MISSOURI Jlinois'

$ws1.Range("A80").Value = 78
$ws1.Range("B80").Value = 87
$ws1.Range("C80").Value = ' JULIAN A. STEYERMARK , COLLECTOR'
$ws1.Range("D80").Value = ' Sept.4 .4 . 1938'
$ws1.Range("E80").Value = ' ( 2 ) Link Pimento along n . Otter Creek Between Tucker and , Jeffry sect . 28 and 29 , 5272 , Q.6 € 3. ) Wayne'
$ws1.Range("F80").Value = ' No. 6626'
$ws1.Range("G80").Value = ' Camptosorus rhizophyllus'
$ws1.Range("H80").Value = ' FLORA OF MISSOURI'

$ws1.Range("A81").Value = 79
$ws1.Range("B81").Value = 88
$ws1.Range("C81").Value = ' JULIAN A. STEYERMARK , COLLECTOR'
$ws1.Range("D81").Value = ' Sept. 4 , 1928'
$ws1.Range("E81").Value = ' Limestone ledges along n . and not - facing Otter Creek and Het fork between Rucker and subway sect . 28 and 29 , 727 , 2.6 , Waynt'
$ws1.Range("F81").Value = ' No. 6626'
$ws1.Range("G81").Value = ' Camptosorus rhizophyllus'
$ws1.Range("H81").Value = ' FLORA OF MISSOURI'
